$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.039.09"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.982.53"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.77"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.56"
$ws.Range("E7").Value = "  +5.02%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.01"
$ws.Range("E12").Value = "  +9.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.28"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.845"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.272.48"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.985.17"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.937.46"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.14"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.17"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.77"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.151"
$ws.Range("E26").Value = "  +8.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.26"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  +16.64%  "
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.87"
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  +5.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.30"
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.35"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.53"
$ws.Range("E39").Value = "  -6.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0967"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.48"
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.28"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.372.56"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.24"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.37"
$ws.Range("E50").Value = "  +6.21%  "
$ws.Range("E51").Value = "  +8.98%  "
